$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText {
    param($range, $text)
    # Prefix with an apostrophe so Excel stores the value as literal text
    # (many of these look like plain numbers, e.g. "574.59") and then
    # clear the resulting quote-prefix cell format so the cell keeps its
    # original (default) style.
    $ws.Range($range).Value = "'" + $text
    $ws.Range($range).ClearFormats()
}

# Row 2 - Bitcoin
Set-PriceText "D2" "62.896.72"
$ws.Range("E2").Value = "  -2.48%  "

# Row 3 - Ethereum
Set-PriceText "D3" "3.412.67"
$ws.Range("E3").Value = "  -2.92%  "

# Row 4 - TetherUSD
Set-PriceText "D4" "0.999"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5 - BNB
Set-PriceText "D5" "574.59"
$ws.Range("E5").Value = "  -2.96%  "

# Row 6 - Solana
Set-PriceText "D6" "127.02"
$ws.Range("E6").Value = "  -5.88%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.02%  "

# Row 8 - LidoStakedEther
Set-PriceText "D8" "3.409.37"
$ws.Range("E8").Value = "  -3.02%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -2.61%  "

# Row 10 - Toncoin
Set-PriceText "D10" "7.41"
$ws.Range("E10").Value = "  -1.56%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -2.78%  "

# Row 12 - Cardano
Set-PriceText "D12" "0.380"
$ws.Range("E12").Value = "  -1.69%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-PriceText "D13" "3.985.42"
$ws.Range("E13").Value = "  -3.08%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -0.77%  "

# Row 15 - was WrappedEther, now ShibaInu
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-PriceText "D15" "0.0000175"
$ws.Range("E15").Value = "  -4.07%  "

# Row 16 - was ShibaInu, now WrappedEther
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-PriceText "D16" "3.407.01"
$ws.Range("E16").Value = "  -2.94%  "

# Row 17 - WrappedBTC
Set-PriceText "D17" "62.940.29"
$ws.Range("E17").Value = "  -2.37%  "

# Row 18 - Avalanche
Set-PriceText "D18" "24.99"
$ws.Range("E18").Value = "  -3.46%  "

# Row 19 - Uniswap
Set-PriceText "D19" "9.63"
$ws.Range("E19").Value = "  -3.00%  "

# Row 20 - Polkadot
Set-PriceText "D20" "5.70"
$ws.Range("E20").Value = "  -1.27%  "

# Row 21 - Chainlink
Set-PriceText "D21" "13.26"
$ws.Range("E21").Value = "  -2.86%  "

# Row 22 - BitcoinCash
Set-PriceText "D22" "378.77"
$ws.Range("E22").Value = "  -3.98%  "

# Row 23 - Polygon
Set-PriceText "D23" "0.562"
$ws.Range("E23").Value = "  -2.59%  "

# Row 24 - WrappedeETH
Set-PriceText "D24" "3.544.21"
$ws.Range("E24").Value = "  -2.97%  "

# Row 25 - Litecoin
Set-PriceText "D25" "72.66"
$ws.Range("E25").Value = "  -2.76%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  -0.12%  "

# Row 27 - PEPE
$ws.Range("E27").Value = "  -8.00%  "

# Row 28 - Binance-PegBSC-USD
$ws.Range("E28").Value = "  -0.01%  "

# Row 29 - RenderToken
Set-PriceText "D29" "7.02"
$ws.Range("E29").Value = "  -5.58%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -4.54%  "

# Row 31 - InternetComputer(DFINITY)
Set-PriceText "D31" "7.91"
$ws.Range("E31").Value = "  -4.84%  "

# Row 32 - was Kaspa, now Fetch.AI
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-PriceText "D32" "1.41"
$ws.Range("E32").Value = "  -4.04%  "

# Row 33 - was Fetch.AI, now Kaspa
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-PriceText "D33" "0.153"
$ws.Range("E33").Value = "  -3.89%  "

# Row 34 - RenzoRestakedETH
Set-PriceText "D34" "3.438.44"
$ws.Range("E34").Value = "  -2.85%  "

# Row 35 - USDe
$ws.Range("E35").Value = "  -0.03%  "

# Row 36 - EthereumClassic
Set-PriceText "D36" "22.95"
$ws.Range("E36").Value = "  -2.28%  "

# Row 37 - NEARProtocol
Set-PriceText "D37" "5.34"
$ws.Range("E37").Value = "  -0.60%  "

# Row 38 - Aptos
Set-PriceText "D38" "6.77"
$ws.Range("E38").Value = "  -2.82%  "

# Row 39 - Monero
Set-PriceText "D39" "164.25"
$ws.Range("E39").Value = "  -1.98%  "

# Row 40 - ImmutableX
Set-PriceText "D40" "1.51"
$ws.Range("E40").Value = "  -3.41%  "

# Row 41 - Hedera
Set-PriceText "D41" "0.0765"
$ws.Range("E41").Value = "  -3.28%  "

# Row 42 - was FirstDigitalUSD, now Mantle
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-PriceText "D42" "0.783"
$ws.Range("E42").Value = "  -3.52%  "

# Row 43 - was Mantle, now FirstDigitalUSD
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-PriceText "D43" "0.999"
$ws.Range("E43").Value = "  -0.13%  "

# Row 44 - OKB
Set-PriceText "D44" "41.62"
$ws.Range("E44").Value = "  -1.65%  "

# Row 45 - Filecoin
Set-PriceText "D45" "4.30"
$ws.Range("E45").Value = "  -3.75%  "

# Row 46 - Stacks
$ws.Range("E46").Value = "  -5.43%  "

# Row 47 - EnergySwap
Set-PriceText "D47" "23.07"
$ws.Range("E47").Value = "  -7.76%  "

# Row 48 - ONDO
$ws.Range("E48").Value = "  -7.35%  "

# Row 49 - Cosmos
Set-PriceText "D49" "6.71"
$ws.Range("E49").Value = "  -1.60%  "

# Row 50 - was Maker, now SuiNetwork
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-PriceText "D50" "0.866"
$ws.Range("E50").Value = "  -4.07%  "

# Row 51 - was SuiNetwork, now Maker
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-PriceText "D51" "2.256.75"
$ws.Range("E51").Value = "  -5.45%  "
